$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 3.7
$ws.Range("I3").Value = 2.3
$ws.Range("L3").Value = 1.53
$ws.Range("M3").Value = 2.38
$ws.Range("N3").Value = 2.7
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 1.57
$ws.Range("Q3").Value = 2.25
$ws.Range("S3").Value = 1.58
$ws.Range("T3").Value = 8
$ws.Range("U3").Value = 17
$ws.Range("X3").Value = 41
$ws.Range("AG3").Value = 10
$ws.Range("AH3").Value = 21

# Row 4
$ws.Range("G4").Value = 2.35
$ws.Range("I4").Value = 2.88
$ws.Range("K4").Value = 9.5
$ws.Range("R4").Value = 1.8
$ws.Range("S4").Value = 1.8
$ws.Range("U4").Value = 11
$ws.Range("X4").Value = 21
$ws.Range("Z4").Value = 9.5
$ws.Range("AE4").Value = 8.5
$ws.Range("AF4").Value = 13
$ws.Range("AH4").Value = 29
$ws.Range("AI4").Value = 23

# Row 5
$ws.Range("R5").Value = 1.5
$ws.Range("S5").Value = 2.37

# Row 6
$ws.Range("G6").Value = 1.75
$ws.Range("H6").Value = 4.1
$ws.Range("I6").Value = 4.2
$ws.Range("R6").Value = 1.67
$ws.Range("AF6").Value = 23

# Row 8
$ws.Range("G8").Value = 1.27
$ws.Range("H8").Value = 4.9
$ws.Range("I8").Value = 10.25
$ws.Range("M8").Value = 3.45
$ws.Range("N8").Value = 1.65
$ws.Range("O8").Value = 1.98
$ws.Range("U8").Value = 5.6
$ws.Range("W8").Value = 7.2
$ws.Range("X8").Value = 11.5
$ws.Range("Z8").Value = 11.75
$ws.Range("AA8").Value = 10.25
$ws.Range("AB8").Value = 27
$ws.Range("AE8").Value = 24
$ws.Range("AF8").Value = 75
$ws.Range("AG8").Value = 32
$ws.Range("AH8").Value = 350
$ws.Range("AJ8").Value = 120

# Row 10
$ws.Range("G10").Value = 1.34
$ws.Range("H10").Value = 4.2
$ws.Range("I10").Value = 9.25
$ws.Range("L10").Value = 1.34
$ws.Range("M10").Value = 2.72
$ws.Range("N10").Value = 2
$ws.Range("O10").Value = 1.65
$ws.Range("P10").Value = 1.4
$ws.Range("Q10").Value = 2.55
$ws.Range("R10").Value = 2.45
$ws.Range("S10").Value = 1.42
$ws.Range("T10").Value = 4.9
$ws.Range("U10").Value = 5.1
$ws.Range("V10").Value = 9.25
$ws.Range("W10").Value = 7.6
$ws.Range("X10").Value = 14
$ws.Range("Y10").Value = 45
$ws.Range("Z10").Value = 8.25
$ws.Range("AA10").Value = 9
$ws.Range("AB10").Value = 30
$ws.Range("AC10").Value = 200
$ws.Range("AD10").Value = 101
$ws.Range("AE10").Value = 18
$ws.Range("AF10").Value = 65
$ws.Range("AG10").Value = 32
$ws.Range("AH10").Value = 300
$ws.Range("AI10").Value = 175
$ws.Range("AJ10").Value = 150

# Row 11
$ws.Range("G11").Value = 2.75
$ws.Range("H11").Value = 2.95
$ws.Range("I11").Value = 2.5
$ws.Range("J11").Value = 1.1
$ws.Range("K11").Value = 5.8
$ws.Range("L11").Value = 1.47
$ws.Range("M11").Value = 2.49
$ws.Range("N11").Value = 2.45
$ws.Range("O11").Value = 1.49
$ws.Range("R11").Value = 2.07
$ws.Range("S11").Value = 1.68
$ws.Range("T11").Value = 5.4
$ws.Range("U11").Value = 10
$ws.Range("V11").Value = 9
$ws.Range("W11").Value = 28
$ws.Range("X11").Value = 25
$ws.Range("Y11").Value = 45
$ws.Range("Z11").Value = 5.2
$ws.Range("AA11").Value = 4.6
$ws.Range("AB11").Value = 14
$ws.Range("AC11").Value = 101
$ws.Range("AD11").Value = 101
$ws.Range("AE11").Value = 5.2
$ws.Range("AF11").Value = 8.800000000000001
$ws.Range("AG11").Value = 8.4
$ws.Range("AH11").Value = 23
$ws.Range("AI11").Value = 22
$ws.Range("AJ11").Value = 40

# Row 12
$ws.Range("I12").Value = 3
$ws.Range("K12").Value = 4.25
$ws.Range("R12").Value = 2.35
$ws.Range("T12").Value = 5.8
$ws.Range("Y12").Value = 60
$ws.Range("Z12").Value = 4.25
$ws.Range("AA12").Value = 5.4
$ws.Range("AE12").Value = 5.8
$ws.Range("AF12").Value = 13
$ws.Range("AG12").Value = 12.5
$ws.Range("AH12").Value = 40

# Row 13
$ws.Range("G13").Value = 2.32
$ws.Range("H13").Value = 2.75
$ws.Range("I13").Value = 3.5
$ws.Range("K13").Value = 4.55
$ws.Range("U13").Value = 9
$ws.Range("V13").Value = 10.75
$ws.Range("Z13").Value = 4.55
$ws.Range("AA13").Value = 5.8
$ws.Range("AB13").Value = 23
$ws.Range("AE13").Value = 6.7
$ws.Range("AF13").Value = 16
$ws.Range("AH13").Value = 55

# Row 14
$ws.Range("S14").Value = 1.63

# Row 17
$ws.Range("J17").Value = 1.07
$ws.Range("K17").Value = 9
$ws.Range("L17").Value = 1.36
$ws.Range("M17").Value = 3
$ws.Range("N17").Value = 2.1
$ws.Range("O17").Value = 1.7
$ws.Range("P17").Value = 1.44
$ws.Range("Q17").Value = 2.63
$ws.Range("R17").Value = 1.95
$ws.Range("S17").Value = 1.8
$ws.Range("V17").Value = 9
$ws.Range("Z17").Value = 8.5
$ws.Range("AB17").Value = 17
$ws.Range("AD17").Value = 401
$ws.Range("AE17").Value = 10
$ws.Range("AI17").Value = 41

# Row 19
$ws.Range("G19").Value = 2.5
$ws.Range("I19").Value = 3.25
$ws.Range("T19").Value = 5.5
$ws.Range("U19").Value = 10
$ws.Range("X19").Value = 26
$ws.Range("AE19").Value = 7

# Row 21
$ws.Range("G21").Value = 3.25
$ws.Range("H21").Value = 3.25
$ws.Range("I21").Value = 2.2
$ws.Range("J21").Value = 1.06
$ws.Range("K21").Value = 10
$ws.Range("T21").Value = 10
$ws.Range("W21").Value = 34
$ws.Range("X21").Value = 26
$ws.Range("AC21").Value = 41
$ws.Range("AE21").Value = 8
$ws.Range("AF21").Value = 11
$ws.Range("AH21").Value = 21

# Row 30
$ws.Range("G30").Value = 3.1
$ws.Range("H30").Value = 3.25
$ws.Range("I30").Value = 2.07
$ws.Range("P30").Value = 1.38
$ws.Range("R30").Value = 1.75
$ws.Range("S30").Value = 1.96
$ws.Range("U30").Value = 13.5
$ws.Range("V30").Value = 9.25
$ws.Range("W30").Value = 32
$ws.Range("X30").Value = 22
$ws.Range("AA30").Value = 5.6
$ws.Range("AB30").Value = 11.5
$ws.Range("AE30").Value = 6.5
$ws.Range("AF30").Value = 8.5
$ws.Range("AG30").Value = 7.4
$ws.Range("AH30").Value = 15.5
$ws.Range("AI30").Value = 13.5

# Row 31
$ws.Range("G31").Value = 3.6
$ws.Range("I31").Value = 1.83
$ws.Range("L31").Value = 1.22
$ws.Range("M31").Value = 3.8
$ws.Range("T31").Value = 10.25
$ws.Range("U31").Value = 17
$ws.Range("V31").Value = 10.25
$ws.Range("W31").Value = 40
$ws.Range("X31").Value = 24
$ws.Range("Y31").Value = 27
$ws.Range("AE31").Value = 7.1
$ws.Range("AF31").Value = 8.25
$ws.Range("AH31").Value = 13
$ws.Range("AI31").Value = 11.25

# Row 34
$ws.Range("G34").Value = 3.1
$ws.Range("I34").Value = 2.22
$ws.Range("L34").Value = 1.44
$ws.Range("M34").Value = 2.4
$ws.Range("N34").Value = 2.27
$ws.Range("P34").Value = 1.5
$ws.Range("Q34").Value = 2.27
$ws.Range("R34").Value = 2
$ws.Range("W34").Value = 40
$ws.Range("Z34").Value = 7.2
$ws.Range("AB34").Value = 18.5
$ws.Range("AE34").Value = 6.1
$ws.Range("AG34").Value = 9.75
$ws.Range("AI34").Value = 22

# Row 36
$ws.Range("G36").Value = 2.9
$ws.Range("I36").Value = 2.35
$ws.Range("K36").Value = 10
$ws.Range("U36").Value = 13
$ws.Range("W36").Value = 29
$ws.Range("X36").Value = 23
$ws.Range("AH36").Value = 23
$ws.Range("AI36").Value = 21

# Row 38
$ws.Range("G38").Value = 3.2
$ws.Range("I38").Value = 1.95
$ws.Range("N38").Value = 1.93
$ws.Range("O38").Value = 1.88
$ws.Range("T38").Value = 10
$ws.Range("U38").Value = 17
$ws.Range("W38").Value = 34
$ws.Range("X38").Value = 26
$ws.Range("AF38").Value = 10

# Row 43
$ws.Range("G43").Value = 2.55
$ws.Range("H43").Value = 3.3
$ws.Range("I43").Value = 2.7
